$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.322190284729004
$ws.Range("B1").Value = 2.561465501785278
$ws.Range("C1").Value = 2.739779233932495
$ws.Range("D1").Value = 3.524868488311768
$ws.Range("E1").Value = 0.8342031240463257
